$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.916.68'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '2.919.29'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.503'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '2.914.82'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.435'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.17'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '3.401.33'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').Value = '62.813.25'
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').Value = '2.914.07'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '432.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.663'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  +5.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.43'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.62%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.269'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('D45').Value = '2.720.23'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0341'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '133.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '356.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('E50').Value = '  +17.10%  '
$ws.Range('E51').Value = '  -0.27%  '
